$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for data rows 2-61.
# Update every occurrence of 45177 to 45178, leaving formatting/style intact.
for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value = 45178
    }
}
